# Auto-applied update of Leve profit calculator columns (H-N) per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 16911.705
$ws.Range("I116").Value = 6714.4287
$ws.Range("J116").Value = 24049.8
$ws.Range("K116").Value = 6714.4287
$ws.Range("L116").Value = 24049.8
$ws.Range("M116").Value = -3272.4287
$ws.Range("N116").Value = -30933.8
$ws.Range("H137").Value = 50005708
$ws.Range("I137").Value = 100003480
$ws.Range("J137").Value = 7935.6
$ws.Range("K137").Value = 300010440
$ws.Range("L137").Value = 23806.8
$ws.Range("M137").Value = -300007890
$ws.Range("N137").Value = -28906.8
$ws.Range("H138").Value = 5050.9473
$ws.Range("J138").Value = 6763.5454
$ws.Range("L138").Value = 20290.6362
$ws.Range("N138").Value = -30570.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12939.645
$ws.Range("I32").Value = 12654.857
$ws.Range("K32").Value = 12654.857
$ws.Range("M32").Value = -12367.857
$ws.Range("H61").Value = 29591400
$ws.Range("I61").Value = 35008180
$ws.Range("J61").Value = 2507499.8
$ws.Range("K61").Value = 35008180
$ws.Range("L61").Value = 2507499.8
$ws.Range("M61").Value = -35007968
$ws.Range("N61").Value = -2507923.8
$ws.Range("H74").Value = 6798
$ws.Range("I74").Value = 1513
$ws.Range("K74").Value = 1513
$ws.Range("M74").Value = -639
$ws.Range("H77").Value = 6798
$ws.Range("I77").Value = 1513
$ws.Range("K77").Value = 7565
$ws.Range("M77").Value = -3197
$ws.Range("H97").Value = 1924.36
$ws.Range("I97").Value = 799.2778
$ws.Range("K97").Value = 799.2778
$ws.Range("M97").Value = -303.2778
$ws.Range("H110").Value = 4688.1304
$ws.Range("I110").Value = 4203.9443
$ws.Range("J110").Value = 6431.2
$ws.Range("K110").Value = 4203.9443
$ws.Range("L110").Value = 6431.2
$ws.Range("M110").Value = -2158.9443
$ws.Range("N110").Value = -10521.2
$ws.Range("H136").Value = 29591400
$ws.Range("I136").Value = 35008180
$ws.Range("J136").Value = 2507499.8
$ws.Range("K136").Value = 105024540
$ws.Range("L136").Value = 7522499.399999999
$ws.Range("M136").Value = -105021990
$ws.Range("N136").Value = -7527599.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1795.0834
$ws.Range("I99").Value = 791.7143
$ws.Range("K99").Value = 791.7143
$ws.Range("M99").Value = 706.2857
$ws.Range("H134").Value = 7144666
$ws.Range("I134").Value = 1948
$ws.Range("K134").Value = 5844
$ws.Range("M134").Value = -3309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 930.875
$ws.Range("I22").Value = 809.4
$ws.Range("J22").Value = 1133.3334
$ws.Range("K22").Value = 809.4
$ws.Range("L22").Value = 1133.3334
$ws.Range("M22").Value = -459.4
$ws.Range("N22").Value = -1833.3334
$ws.Range("H31").Value = 25003304
$ws.Range("I31").Value = 33336170
$ws.Range("K31").Value = 33336170
$ws.Range("M31").Value = -33335875
$ws.Range("H34").Value = 25003304
$ws.Range("I34").Value = 33336170
$ws.Range("K34").Value = 33336170
$ws.Range("M34").Value = -33335968
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 30000
$ws.Range("I51").Value = 20000
$ws.Range("K51").Value = 20000
$ws.Range("M51").Value = -19264
$ws.Range("H61").Value = 30000
$ws.Range("I61").Value = 20000
$ws.Range("K61").Value = 20000
$ws.Range("M61").Value = -19652

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52634130
$ws.Range("I12").Value = 166669540
$ws.Range("J12").Value = 2397.6155
$ws.Range("K12").Value = 500008620
$ws.Range("L12").Value = 7192.8465
$ws.Range("M12").Value = -500008447
$ws.Range("N12").Value = -7538.8465
$ws.Range("H87").Value = 1439
$ws.Range("J87").Value = 1158
$ws.Range("L87").Value = 3474
$ws.Range("N87").Value = -5970
$ws.Range("H90").Value = 1439
$ws.Range("J90").Value = 1158
$ws.Range("L90").Value = 10422
$ws.Range("N90").Value = -22902
$ws.Range("H113").Value = 1970.0625
$ws.Range("I113").Value = 1730.375
$ws.Range("J113").Value = 2209.75
$ws.Range("K113").Value = 5191.125
$ws.Range("L113").Value = 6629.25
$ws.Range("M113").Value = -3021.125
$ws.Range("N113").Value = -10969.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13203580
$ws.Range("I22").Value = 33003976
$ws.Range("J22").Value = 3316.1667
$ws.Range("K22").Value = 33003976
$ws.Range("L22").Value = 3316.1667
$ws.Range("M22").Value = -33003681
$ws.Range("N22").Value = -3906.1667
$ws.Range("H27").Value = 13203580
$ws.Range("I27").Value = 33003976
$ws.Range("J27").Value = 3316.1667
$ws.Range("K27").Value = 33003976
$ws.Range("L27").Value = 3316.1667
$ws.Range("M27").Value = -33003869
$ws.Range("N27").Value = -3530.1667
$ws.Range("H43").Value = 8000
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7807
$ws.Range("H80").Value = 33987.5
$ws.Range("J80").Value = 33987.5
$ws.Range("L80").Value = 33987.5
$ws.Range("N80").Value = -36233.5
$ws.Range("H83").Value = 33987.5
$ws.Range("J83").Value = 33987.5
$ws.Range("L83").Value = 101962.5
$ws.Range("N83").Value = -113194.5
$ws.Range("H132").Value = 3045.5305
$ws.Range("I132").Value = 2397.361
$ws.Range("J132").Value = 4840.4614
$ws.Range("K132").Value = 7192.083
$ws.Range("L132").Value = 14521.3842
$ws.Range("M132").Value = -4662.083
$ws.Range("N132").Value = -19581.3842
$ws.Range("H136").Value = 3501.9744
$ws.Range("I136").Value = 3256.0645
$ws.Range("J136").Value = 4454.875
$ws.Range("K136").Value = 9768.193499999999
$ws.Range("L136").Value = 13364.625
$ws.Range("M136").Value = -7218.193499999999
$ws.Range("N136").Value = -18464.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 1655.7142
$ws.Range("I100").Value = 870.9091
$ws.Range("K100").Value = 1741.8182
$ws.Range("M100").Value = -1200.8182
$ws.Range("H126").Value = 3998.0715
$ws.Range("I126").Value = 4498.1665
$ws.Range("K126").Value = 13494.4995
$ws.Range("M126").Value = -11024.4995
$ws.Range("H132").Value = 279214.66
$ws.Range("I132").Value = 1150.1305
$ws.Range("J132").Value = 771175
$ws.Range("K132").Value = 3450.3915
$ws.Range("L132").Value = 2313525
$ws.Range("M132").Value = -920.3914999999997
$ws.Range("N132").Value = -2318585
$ws.Range("H136").Value = 275667.66
$ws.Range("J136").Value = 1252686.9
$ws.Range("L136").Value = 3758060.7
$ws.Range("N136").Value = -3763160.7
